# Insurance_Payments.xlsx — "Add files via upload" re-edit.
#
# The author re-typed the E248:E281 "avg payment %" column: every value in
# that range ends up exactly 100x its previous value (e.g. 0.122 -> 12.2,
# 1E-3 -> 0.1, 0.272 -> 27.2, ...). The cells keep their existing percentage
# number format (style id 13 / numFmtId 10), so this reads as the author
# re-entering the percentages as whole numbers instead of decimals.
#
# The view was also left scrolled to a different spot with a different
# active cell/selection (F247, which is where the long-standing review
# comment lives) instead of the prior A248:F281 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Insurance Payment Avgs")

for ($r = 248; $r -le 281; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val * 100
    }
}

# Update the view state to match: scrolled down a bit further, with F247
# (the commented cell) selected instead of the old A248:F281 block.
$ws.Activate() | Out-Null
$ws.Range("F247").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 229
$excel.ActiveWindow.ScrollColumn = 1
